# Update "想去人数" (number of people interested) figures in the
# "展览" and "全部类型" worksheets to the refreshed values from the
# regenerated gh-pages data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 5713
$wsExpo.Range("F3").Value = 84
$wsExpo.Range("F4").Value = 9
$wsExpo.Range("F5").Value = 961
$wsExpo.Range("F6").Value = 158
$wsExpo.Range("F7").Value = 2609
$wsExpo.Range("F9").Value = 185
$wsExpo.Range("F11").Value = 97
$wsExpo.Range("F12").Value = 38
$wsExpo.Range("F13").Value = 2453
$wsExpo.Range("F14").Value = 495

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5713
$wsAll.Range("F3").Value = 84
$wsAll.Range("F4").Value = 9
$wsAll.Range("F6").Value = 961
$wsAll.Range("F7").Value = 158
$wsAll.Range("F8").Value = 2609
$wsAll.Range("F10").Value = 185
$wsAll.Range("F13").Value = 97
$wsAll.Range("F14").Value = 38
$wsAll.Range("F15").Value = 2453
$wsAll.Range("F16").Value = 495
